$wb = $excel.ActiveWorkbook

# Metadata sheet updates
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/ValueSet/insurance-plan-type"
$wsMeta.Range("B3").Value = "8.0.0"
$wsMeta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$wsMeta.Range("B9").Value = "LinuxForHealth Team"

# "Include from Insurance Plan T" sheet update
$wsCode = $wb.Worksheets.Item("Include from Insurance Plan T")
$wsCode.Range("B4").Value = "http://linuxforhealth.org/fhir/cdm/CodeSystem/insurance-plan-type"
